# rerun lm; add permuted models brs, sdnn
# Update regression-output values (est, l95, u95, p, interaction) for rows 2-21
# in columns C:G of the active worksheet, as produced by the rerun of the
# linear-model analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(-0.44, -3.16, 2.28, 0.7496, 0.07374),
    @(3.71, 0.93, 6.5, 0.00931, 0.36459),
    @(-1.04, -3.74, 1.67, 0.45055, 0.02783),
    @(3.83, 1.16, 6.5, 0.00515, 0.2306),
    @(-1.81, -5.19, 1.57, 0.2924, 0.68337),
    @(4.06, 1.46, 6.67, 0.00241, 0.73356),
    @(0.57, -2.15, 3.29, 0.67779, 0.03033),
    @(-0.33, -3.07, 2.42, 0.81381, 0.03644),
    @(3.31, -0.12, 6.74, 0.05859, 0.55239),
    @(3.91, 0.88, 6.94, 0.01182, 0.16821),
    @(2.92, 0.29, 5.55, 0.03003, 0.07374),
    @(5.18, 2.09, 8.27, 0.00116, 0.36459),
    @(3.7, 1.04, 6.36, 0.00665, 0.02783),
    @(4.71, 1.92, 7.51, 0.00106, 0.2306),
    @(0.11, -3.18, 3.4, 0.94865, 0.68337),
    @(3.1, 0.46, 5.74, 0.02178, 0.73356),
    @(4.43, 1.84, 7.02, 0.00093, 0.03033),
    @(3.73, 1.06, 6.4, 0.00642, 0.03644),
    @(4.4, 1.25, 7.56, 0.00648, 0.55239),
    @(2.56, -0.58, 5.7, 0.10987, 0.16821)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $vals = $newValues[$i]
    $ws.Cells.Item($row, 3).Value = $vals[0]   # C - est
    $ws.Cells.Item($row, 4).Value = $vals[1]   # D - l95
    $ws.Cells.Item($row, 5).Value = $vals[2]   # E - u95
    $ws.Cells.Item($row, 6).Value = $vals[3]   # F - p
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - interaction
}
